# Slide 25 ("Two-stage Compiler") contains several copies of a small
# diagram built from nested groups. Two of the little text boxes inside
# one copy of that diagram read "X86/A" (capital X) and need to become
# "x86/A" (lower-case x), with a slightly narrower bounding box to match.
#
# This headless COM runtime's Shape.Left/Top/Width/Height setters write
# straight into the shape's own <a:off>/<a:ext> (the group-local child
# coordinate system) using a float32 points<->EMU round trip, instead of
# inverse-transforming through the parent group's chOff/chExt like real
# PowerPoint does. Because of that float32 rounding the naive
# `emu/12700.0` conversion frequently lands one EMU below the intended
# integer, so a tiny epsilon is added before converting back to points to
# make sure the stored integer lands exactly on target. Editing the text
# of an autosize (spAutoFit) text box also recalculates Height in real
# EMU, so Height is restored (to its unchanged value) right after any
# text edit.

function Get-ShapeById($shapes, $targetId) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $sh = $shapes.Item($i)
        if ($sh.Id -eq $targetId) {
            return $sh
        }
    }
    return $null
}

function Set-ShapeEmu($shape, $leftEmu, $topEmu, $widthEmu, $heightEmu) {
    # Small epsilon compensates float32 precision loss so that
    # round-tripping emu -> points -> emu lands back on the exact integer.
    $epsilon = 0.00004
    if ($topEmu -ne $null)    { $shape.Top    = ($topEmu    + $epsilon) / 12700.0 }
    if ($heightEmu -ne $null) { $shape.Height = ($heightEmu + $epsilon) / 12700.0 }
    if ($leftEmu -ne $null)   { $shape.Left   = ($leftEmu   + $epsilon) / 12700.0 }
    if ($widthEmu -ne $null)  { $shape.Width  = ($widthEmu  + $epsilon) / 12700.0 }
}

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(25)

# Top-level group ("Group 64", shape id 65) holding the copy of the
# diagram whose two labels change in this edit.
$topGroup = Get-ShapeById $s.Shapes 65

# --- Shape id 96 ("Text Box 82"): single run "X86/A" -> "x86/A" --------
$shape96 = Get-ShapeById $topGroup.GroupItems 96
$shape96.TextFrame.TextRange.Text = "x86/A"
# off x="1601"->"1611", ext cx="447"->"427"; y/cy (3479/204) unchanged.
Set-ShapeEmu $shape96 1611 3479 427 204

# --- Shape id 86 ("Text Box 104"): first run "X86/A " -> "x86/A " ------
# (keeps its second run " x86", with the Symbol-font formatting, intact)
$shape86 = Get-ShapeById $topGroup.GroupItems 86
$firstRun = $shape86.TextFrame.TextRange.Characters(1, 6)
$firstRun.Text = "x86/A "
# off x="672"->"682", ext cx="822"->"802"; y/cy (2557/204) unchanged.
Set-ShapeEmu $shape86 682 2557 802 204
